$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2285.9834
$ws.Range("J17").Value = 2285.9834
$ws.Range("L17").Value = 6857.9502
$ws.Range("N17").Value = -7193.9502
$ws.Range("K51").Value = 166666670
$ws.Range("N51").Value = -14218
$ws.Range("I51").Value = 166666670
$ws.Range("H51").Value = 33343934
$ws.Range("J51").Value = 13250
$ws.Range("L51").Value = 13250
$ws.Range("M51").Value = -166666186
$ws.Range("I54").Value = 19663.334
$ws.Range("H54").Value = 27268.5
$ws.Range("J54").Value = 50084
$ws.Range("M54").Value = -19177.334
$ws.Range("L54").Value = 50084
$ws.Range("K54").Value = 19663.334
$ws.Range("N54").Value = -51056
$ws.Range("L62").Value = 2200
$ws.Range("K62").Value = 4717.6
$ws.Range("M62").Value = -4093.6
$ws.Range("N62").Value = -3448
$ws.Range("I62").Value = 4717.6
$ws.Range("H62").Value = 4187.579
$ws.Range("J62").Value = 2200
$ws.Range("N65").Value = -17240
$ws.Range("J65").Value = 2200
$ws.Range("I65").Value = 4717.6
$ws.Range("L65").Value = 11000
$ws.Range("H65").Value = 4187.579
$ws.Range("K65").Value = 23588
$ws.Range("M65").Value = -20468
$ws.Range("I74").Value = 3795.3845
$ws.Range("H74").Value = 3755.2942
$ws.Range("J74").Value = 3625
$ws.Range("M74").Value = -2859.3845
$ws.Range("L74").Value = 3625
$ws.Range("K74").Value = 3795.3845
$ws.Range("N74").Value = -5497
$ws.Range("N77").Value = -27485
$ws.Range("I77").Value = 3795.3845
$ws.Range("H77").Value = 3755.2942
$ws.Range("J77").Value = 3625
$ws.Range("L77").Value = 18125
$ws.Range("M77").Value = -14296.9225
$ws.Range("K77").Value = 18976.9225
$ws.Range("L92").Value = 3459.1667
$ws.Range("M92").Value = 778
$ws.Range("N92").Value = -5955.1667
$ws.Range("J92").Value = 3459.1667
$ws.Range("I92").Value = 470
$ws.Range("H92").Value = 1249.7826
$ws.Range("K92").Value = 470
$ws.Range("N96").Value = -600090046
$ws.Range("I96").Value = 2326.1765
$ws.Range("H96").Value = 45462956
$ws.Range("J96").Value = 200029100
$ws.Range("L96").Value = 600087300
$ws.Range("K96").Value = 6978.529500000001
$ws.Range("M96").Value = -5605.529500000001
$ws.Range("I137").Value = 5495491
$ws.Range("H137").Value = 3210557.5
$ws.Range("K137").Value = 16486473
$ws.Range("M137").Value = -16483923
$ws.Range("L137").Value = 34950.60000000001
$ws.Range("N137").Value = -40050.60000000001
$ws.Range("J137").Value = 11650.2
$ws.Range("N138").Value = -25872.1061
$ws.Range("I138").Value = 1447.1852
$ws.Range("H138").Value = 2996.1738
$ws.Range("J138").Value = 5197.3687
$ws.Range("K138").Value = 4341.5556
$ws.Range("L138").Value = 15592.1061
$ws.Range("M138").Value = 798.4444000000003

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("I32").Value = 2995.9707
$ws.Range("H32").Value = 5966.436
$ws.Range("K32").Value = 2995.9707
$ws.Range("M32").Value = -2708.9707
$ws.Range("I53").Value = 10750
$ws.Range("H53").Value = 10750
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("M53").Value = -10068
$ws.Range("K53").Value = 10750
$ws.Range("N53").ClearContents()
$ws.Range("H54").Value = 8664.166999999999
$ws.Range("J54").Value = 8664.166999999999
$ws.Range("L54").Value = 8664.166999999999
$ws.Range("N54").Value = -10202.167
$ws.Range("M97").Value = -1710.6667
$ws.Range("I97").Value = 2206.6667
$ws.Range("H97").Value = 2321.8333
$ws.Range("K97").Value = 2206.6667
$ws.Range("I122").Value = 1363.6875
$ws.Range("H122").Value = 1452.4048
$ws.Range("K122").Value = 4091.0625
$ws.Range("M122").Value = -1641.0625
$ws.Range("H137").Value = 35000
$ws.Range("L137").Value = 35000
$ws.Range("N137").Value = -45200
$ws.Range("J137").Value = 35000

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("K94").Value = 2313.7144
$ws.Range("N94").Value = -2625.3334
$ws.Range("L94").Value = 1723.3334
$ws.Range("I94").Value = 2313.7144
$ws.Range("H94").Value = 2136.6
$ws.Range("J94").Value = 1723.3334
$ws.Range("M94").Value = -1862.7144
$ws.Range("H99").Value = 2086.5117
$ws.Range("J99").Value = 1875
$ws.Range("K99").Value = 2168.3872
$ws.Range("L99").Value = 1875
$ws.Range("M99").Value = -670.3872000000001
$ws.Range("N99").Value = -4871
$ws.Range("I99").Value = 2168.3872
$ws.Range("N110").ClearContents()
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I7").Value = 0
$ws.Range("H7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("H16").Value = 1790.5
$ws.Range("J16").Value = 1682.6
$ws.Range("L16").Value = 1682.6
$ws.Range("M16").Value = -1683.3334
$ws.Range("K16").Value = 1970.3334
$ws.Range("N16").Value = -2256.6
$ws.Range("I16").Value = 1970.3334
$ws.Range("M31").Value = -1160.0571
$ws.Range("N31").Value = -14503119
$ws.Range("J31").Value = 14502529
$ws.Range("I31").Value = 1455.0571
$ws.Range("L31").Value = 14502529
$ws.Range("H31").Value = 5751881
$ws.Range("K31").Value = 1455.0571
$ws.Range("I34").Value = 1455.0571
$ws.Range("H34").Value = 5751881
$ws.Range("K34").Value = 1455.0571
$ws.Range("J34").Value = 14502529
$ws.Range("L34").Value = 14502529
$ws.Range("M34").Value = -1253.0571
$ws.Range("N34").Value = -14502933
$ws.Range("N113").Value = -6022.6
$ws.Range("I113").Value = 1970.3334
$ws.Range("H113").Value = 1790.5
$ws.Range("J113").Value = 1682.6
$ws.Range("K113").Value = 1970.3334
$ws.Range("L113").Value = 1682.6
$ws.Range("M113").Value = 199.6666
$ws.Range("L116").Value = 43528
$ws.Range("N116").Value = -52706
$ws.Range("H116").Value = 43528
$ws.Range("J116").Value = 43528
$ws.Range("K132").Value = 4720.2858
$ws.Range("M132").Value = -2190.2858
$ws.Range("I132").Value = 1573.4286
$ws.Range("H132").Value = 468697.9

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I5").Value = 8894
$ws.Range("H5").Value = 5179.4346
$ws.Range("K5").Value = 26682
$ws.Range("M5").Value = -26570
$ws.Range("N122").Value = -74610.40000000001
$ws.Range("I122").Value = 430.11765
$ws.Range("L122").Value = 69710.40000000001
$ws.Range("H122").Value = 3859.25
$ws.Range("J122").Value = 7745.6
$ws.Range("K122").Value = 3871.05885
$ws.Range("M122").Value = -1421.05885
$ws.Range("L132").Value = 22081.5
$ws.Range("K132").Value = 12461.1426
$ws.Range("M132").Value = -9931.142600000001
$ws.Range("N132").Value = -27141.5
$ws.Range("I132").Value = 1384.5714
$ws.Range("H132").Value = 1919.0358
$ws.Range("J132").Value = 2453.5
$ws.Range("I135").Value = 8894
$ws.Range("H135").Value = 5179.4346
$ws.Range("M135").Value = -77511
$ws.Range("K135").Value = 80046

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("M97").Value = -2424.625
$ws.Range("N97").Value = -11666
$ws.Range("I97").Value = 2920.625
$ws.Range("H97").Value = 4144.8423
$ws.Range("J97").Value = 10674
$ws.Range("L97").Value = 10674
$ws.Range("K97").Value = 2920.625
$ws.Range("N113").Value = -6473.3333
$ws.Range("I113").Value = 2101.2222
$ws.Range("H113").Value = 2109.25
$ws.Range("J113").Value = 2133.3333
$ws.Range("K113").Value = 2101.2222
$ws.Range("L113").Value = 2133.3333
$ws.Range("M113").Value = 68.77779999999984
$ws.Range("N122").Value = -8300.0002
$ws.Range("I122").Value = 1200.5714
$ws.Range("L122").Value = 3400.0002
$ws.Range("H122").Value = 1169.5385
$ws.Range("J122").Value = 1133.3334
$ws.Range("K122").Value = 3601.7142
$ws.Range("M122").Value = -1151.7142

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("K51").Value = 20000
$ws.Range("N51").Value = -30951
$ws.Range("I51").Value = 20000
$ws.Range("H51").Value = 26663.334
$ws.Range("J51").Value = 29995
$ws.Range("M51").Value = -19522
$ws.Range("L51").Value = 29995
$ws.Range("H93").Value = 1114.3158
$ws.Range("J93").Value = 1148.1428
$ws.Range("K93").Value = 1019.6
$ws.Range("L93").Value = 1148.1428
$ws.Range("M93").Value = 228.4
$ws.Range("N93").Value = -3644.1428
$ws.Range("I93").Value = 1019.6

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I58").Value = 0
$ws.Range("H58").Value = 30044.666
$ws.Range("J58").Value = 30044.666
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("L58").Value = 30044.666
$ws.Range("N58").Value = -30660.666
$ws.Range("J100").Value = 492.25
$ws.Range("I100").Value = 510.93332
$ws.Range("L100").Value = 984.5
$ws.Range("K100").Value = 1021.86664
$ws.Range("M100").Value = -480.86664
$ws.Range("N100").Value = -2066.5
$ws.Range("H100").Value = 507
$ws.Range("L132").Value = 9910.799999999999
$ws.Range("K132").Value = 4663400.4
$ws.Range("M132").Value = -4660870.4
$ws.Range("N132").Value = -14970.8
$ws.Range("I132").Value = 1554466.8
$ws.Range("H132").Value = 1146266
$ws.Range("J132").Value = 3303.6

Write-Output "Applied 237 cell updates across 8 sheets"